$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1346153846153846
$ws.Range("C2").Value = 0.6875
$ws.Range("J2").Value = 0.004807692307692308
$ws.Range("P2").Value = 0.1153846153846154
$ws.Range("S2").Value = 0.0576923076923077
$ws.Range("B3").Value = 0.02649006622516556
$ws.Range("C3").Value = 0.05298013245033113
$ws.Range("J3").Value = 0.01324503311258278
$ws.Range("P3").Value = 0.7748344370860927
$ws.Range("S3").Value = 0.1324503311258278
$ws.Range("J4").Value = 0.03333333333333333
$ws.Range("P4").Value = 0.7666666666666667
$ws.Range("S4").Value = 0.2
$ws.Range("B6").Value = 0.08571428571428572
$ws.Range("D6").Value = 0.01714285714285714
$ws.Range("F6").Value = 0.06857142857142857
$ws.Range("J6").Value = 0.1542857142857143
$ws.Range("O6").Value = 0.01142857142857143
$ws.Range("Q6").Value = 0.1314285714285714
$ws.Range("R6").Value = 0.12
$ws.Range("S6").Value = 0.4114285714285714
$ws.Range("B7").Value = 0.08552631578947369
$ws.Range("D7").Value = 0.0131578947368421
$ws.Range("E7").Value = 0.006578947368421052
$ws.Range("F7").Value = 0.05263157894736842
$ws.Range("J7").Value = 0.131578947368421
$ws.Range("O7").Value = 0.0131578947368421
$ws.Range("Q7").Value = 0.1447368421052632
$ws.Range("R7").Value = 0.1052631578947368
$ws.Range("S7").Value = 0.4473684210526316
$ws.Range("B8").Value = 0.05084745762711865
$ws.Range("D8").Value = 0.02259887005649718
$ws.Range("E8").Value = 0.002824858757062147
$ws.Range("F8").Value = 0.05932203389830509
$ws.Range("J8").Value = 0.1045197740112994
$ws.Range("O8").Value = 0.008474576271186441
$ws.Range("Q8").Value = 0.1864406779661017
$ws.Range("R8").Value = 0.1129943502824859
$ws.Range("S8").Value = 0.4519774011299435
$ws.Range("B9").Value = 0.09268292682926829
$ws.Range("D9").Value = 0.01951219512195122
$ws.Range("E9").Value = 0.004878048780487805
$ws.Range("F9").Value = 0.04878048780487805
$ws.Range("J9").Value = 0.1317073170731707
$ws.Range("O9").Value = 0.01463414634146342
$ws.Range("Q9").Value = 0.1853658536585366
$ws.Range("R9").Value = 0.1170731707317073
$ws.Range("S9").Value = 0.3853658536585366
$ws.Range("B10").Value = 0.0944055944055944
$ws.Range("D10").Value = 0.01311188811188811
$ws.Range("E10").Value = 0.0008741258741258741
$ws.Range("F10").Value = 0.07342657342657342
$ws.Range("J10").Value = 0.1031468531468532
$ws.Range("O10").Value = 0.01660839160839161
$ws.Range("Q10").Value = 0.1853146853146853
$ws.Range("R10").Value = 0.1215034965034965
$ws.Range("S10").Value = 0.3916083916083916
$ws.Range("G11").Value = 0.1528384279475982
$ws.Range("J11").Value = 0.1004366812227074
$ws.Range("K11").Value = 0.2096069868995633
$ws.Range("L11").Value = 0.5240174672489083
$ws.Range("S11").Value = 0.01310043668122271
$ws.Range("G12").Value = 0.7559055118110236
$ws.Range("J12").Value = 0.1417322834645669
$ws.Range("K12").Value = 0.01574803149606299
$ws.Range("L12").Value = 0.06299212598425197
$ws.Range("S12").Value = 0.02362204724409449
$ws.Range("F13").Value = 0.03333333333333333
$ws.Range("G13").Value = 0.8333333333333334
$ws.Range("J13").Value = 0.1333333333333333
$ws.Range("F15").Value = 0.01282051282051282
$ws.Range("H15").Value = 0.1111111111111111
$ws.Range("I15").Value = 0.06837606837606838
$ws.Range("J15").Value = 0.4017094017094017
$ws.Range("K15").Value = 0.06837606837606838
$ws.Range("M15").Value = 0.004273504273504274
$ws.Range("O15").Value = 0.08974358974358974
$ws.Range("S15").Value = 0.2435897435897436
$ws.Range("H16").Value = 0.15527950310559
$ws.Range("I16").Value = 0.08074534161490683
$ws.Range("J16").Value = 0.4906832298136646
$ws.Range("K16").Value = 0.08695652173913043
$ws.Range("M16").Value = 0.02484472049689441
$ws.Range("O16").Value = 0.09937888198757763
$ws.Range("S16").Value = 0.06211180124223602
$ws.Range("F17").Value = 0.01111111111111111
$ws.Range("H17").Value = 0.1666666666666667
$ws.Range("I17").Value = 0.1138888888888889
$ws.Range("J17").Value = 0.4666666666666667
$ws.Range("K17").Value = 0.06388888888888888
$ws.Range("M17").Value = 0.008333333333333333
$ws.Range("O17").Value = 0.06111111111111111
$ws.Range("S17").Value = 0.1083333333333333
$ws.Range("F18").Value = 0.01244813278008299
$ws.Range("H18").Value = 0.1618257261410788
$ws.Range("I18").Value = 0.0954356846473029
$ws.Range("J18").Value = 0.4896265560165975
$ws.Range("K18").Value = 0.07053941908713693
$ws.Range("M18").Value = 0.02489626556016597
$ws.Range("O18").Value = 0.07883817427385892
$ws.Range("S18").Value = 0.06639004149377593
$ws.Range("F19").Value = 0.01291512915129151
$ws.Range("H19").Value = 0.1881918819188192
$ws.Range("I19").Value = 0.1033210332103321
$ws.Range("J19").Value = 0.3902214022140221
$ws.Range("K19").Value = 0.0959409594095941
$ws.Range("M19").Value = 0.01476014760147601
$ws.Range("N19").Value = 0.001845018450184502
$ws.Range("O19").Value = 0.08763837638376384
$ws.Range("S19").Value = 0.1051660516605166
